# Insert 5 new "Patient" demographic fields (age, sex, weight, ethnicity, race)
# into the device_event_fields reference sheet, right after the existing
# "patient.patient_sequence_number" row (row 52), and before
# "patient.patient_problems" (row 53).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device_event_fields")

# Insert 5 blank rows before current row 53 ("patient.patient_problems"),
# shifting existing rows 53+ downward.
$insertRange = $ws.Range("A53:D57")
$insertRange.EntireRow.Insert()

$newRows = @(
    @{ A = "Patient"; B = "patient.patient_age";       C = "string"; D = "Patient's age. This is an .exact field. It has been indexed both as its exact string content, and also tokenized." },
    @{ A = "Patient"; B = "patient.patient_sex";        C = "string"; D = "Patient's gender. This is an .exact field. It has been indexed both as its exact string content, and also tokenized." },
    @{ A = "Patient"; B = "patient.patient_weight";     C = "string"; D = "Patient's weight. This is an .exact field. It has been indexed both as its exact string content, and also tokenized." },
    @{ A = "Patient"; B = "patient.patient_ethnicity";  C = "string"; D = "Patient's ethnicity. This is an .exact field. It has been indexed both as its exact string content, and also tokenized." },
    @{ A = "Patient"; B = "patient.patient_race";       C = "string"; D = "Patient's race. This is an .exact field. It has been indexed both as its exact string content, and also tokenized." }
)

# Populate column by column (A, then B, then C, then D) across all five new
# rows so that new shared-string entries are appended in the same order as
# the source workbook (all "B" values first, then all "D" values).
$row = 53
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 1).Value = $rec.A
    $row++
}

$row = 53
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 2).Value = $rec.B
    $row++
}

$row = 53
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 3).Value = $rec.C
    $row++
}

$row = 53
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 4).Value = $rec.D
    $row++
}

$wb.Save()
